$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric to Excel; force them to remain text
# by applying a text number format before assignment, then resetting the
# style back to Normal so no stray style survives on the cell.
$textForceCells = @("D5", "D6", "D10", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D36", "D39", "D40", "D41", "D43", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.791.83"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "3.133.30"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "602.18"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("D6").Value = "142.99"
$ws.Range("E6").Value = "  -3.53%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.128.66"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").Value = "0.149"
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("D11").Value = "5.38"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").Value = "0.466"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").Value = "35.04"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").Value = "3.639.99"
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").Value = "63.815.30"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "3.132.90"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").Value = "6.81"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").Value = "484.99"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").Value = "14.67"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "0.707"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").Value = "7.64"
$ws.Range("E23").Value = "  -4.88%  "
$ws.Range("D24").Value = "86.76"
$ws.Range("E24").Value = "  +3.38%  "
$ws.Range("D25").Value = "13.43"
$ws.Range("E25").Value = "  -1.92%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "2.75"
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("D28").Value = "8.26"
$ws.Range("E28").Value = "  -3.60%  "
$ws.Range("D29").Value = "7.00"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("D30").Value = "2.06"
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("D31").Value = "27.10"
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  -8.32%  "
$ws.Range("E34").Value = "  -3.03%  "
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("D36").Value = "5.99"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").Value = "0.0₃0744"
$ws.Range("E38").Value = "  -5.64%  "
$ws.Range("D39").Value = "2.96"
$ws.Range("E39").Value = "  -8.50%  "
$ws.Range("D40").Value = "437.66"
$ws.Range("E40").Value = "  -4.73%  "
$ws.Range("D41").Value = "0.0394"
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "8.27"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").Value = "2.873.18"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D45").Value = "0.260"
$ws.Range("E45").Value = "  -3.43%  "
$ws.Range("D46").Value = "2.20"
$ws.Range("E46").Value = "  -5.74%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "2.37"
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").Value = "25.88"
$ws.Range("E49").Value = "  -2.60%  "
$ws.Range("D50").Value = "0.114"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("E51").Value = "  +0.97%  "

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
